$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("E5").Value = "2016-03-24 04:41:20"
